$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing 2021 column (O) into the new 2022 column (P)
# for the header/data rows (3-5), matching the existing per-row styles.
$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial(-4122)

$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)

$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Set the new values
$ws.Range("P3").Value = 2022
$ws.Range("P4").Value = 15
$ws.Range("P5").Value = 2130.4

# Update selection to match the committed workbook state
$ws.Range("P6").Select()
